$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing data (B:M) to (C:N)
$ws.Range("B1").EntireColumn.Insert()

# New header for inserted column
$ws.Range("B1").Value = "ID"

# New building-id values for the inserted column
$ws.Range("B2").Value = "building_towncenter"
$ws.Range("B3").Value = "building_farm"
$ws.Range("B4").Value = "building_factory"
$ws.Range("B5").Value = "building_filterationplant"
$ws.Range("B6").Value = "building_house"

# Update selection to match the authored workbook state
$ws.Range("B2:B6").Select()
